$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Marking" row -> Right column: 3 -> 5
$ws.Range("B11").Value = 5

# "Total" row -> Right column: 84 -> 140
$ws.Range("B12").Value = 140

# "Total" row -> Max column label: "84/84" -> "140/140"
$ws.Range("E12").Value = "140/140"
